$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '26.092.91'
$ws.Range('E2').Value = '  -0.74%  '
$ws.Range('D3').Value = '1.651.11'
$ws.Range('E3').Value = '  -0.83%  '
$ws.Range('E4').Value = '  -0.44%  '
Set-TextValue $ws.Range('D5') '218.07'
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('E6').Value = '  -2.23%  '
$ws.Range('E7').Value = '  -0.42%  '
Set-TextValue $ws.Range('D9') '0.06301'
$ws.Range('E9').Value = '  -1.24%  '
Set-TextValue $ws.Range('D10') '20.54'
$ws.Range('E10').Value = '  -0.18%  '
Set-TextValue $ws.Range('D11') '0.07811'
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('E12').Value = '  -1.75%  '
$ws.Range('D13').Value = '1.651.78'
$ws.Range('E13').Value = '  -0.68%  '
$ws.Range('D14').Value = '1.878.17'
$ws.Range('E14').Value = '  -0.77%  '
$ws.Range('E15').Value = '  +0.52%  '
$ws.Range('D16').Value = '0.0₅8017'
$ws.Range('E16').Value = '  -2.39%  '
$ws.Range('E17').Value = '  -1.12%  '
$ws.Range('D18').Value = '26.078.24'
$ws.Range('E18').Value = '  -0.84%  '
$ws.Range('E19').Value = '  -0.40%  '
Set-TextValue $ws.Range('D20') '4.642'
Set-TextValue $ws.Range('D21') '194.72'
$ws.Range('E21').Value = '  +0.63%  '
Set-TextValue $ws.Range('D22') '10.08'
$ws.Range('E22').Value = '  -1.21%  '
Set-TextValue $ws.Range('D23') '5.953'
$ws.Range('E23').Value = '  -1.37%  '
Set-TextValue $ws.Range('D24') '1.007'
$ws.Range('E24').Value = '  -0.41%  '
Set-TextValue $ws.Range('D25') '146.86'
$ws.Range('E25').Value = '  +0.88%  '
$ws.Range('E26').Value = '  -1.96%  '
$ws.Range('E27').Value = '  -0.21%  '
$ws.Range('E28').Value = '  -1.32%  '
Set-TextValue $ws.Range('D29') '1.472'
$ws.Range('E29').Value = '  -0.60%  '
Set-TextValue $ws.Range('D30') '0.05685'
$ws.Range('E30').Value = '  -3.16%  '
$ws.Range('E31').Value = '  -1.07%  '
Set-TextValue $ws.Range('D32') '3.481'
$ws.Range('E32').Value = '  -3.50%  '
Set-TextValue $ws.Range('D33') '3.364'
$ws.Range('E33').Value = '  +2.61%  '
Set-TextValue $ws.Range('D34') '1.593'
$ws.Range('E34').Value = '  -0.99%  '
Set-TextValue $ws.Range('D35') '2.801'
$ws.Range('E35').Value = '  -0.91%  '
Set-TextValue $ws.Range('D36') '0.9503'
$ws.Range('E36').Value = '  -1.33%  '
$ws.Range('E37').Value = '  -0.21%  '
Set-TextValue $ws.Range('D38') '0.5660'
$ws.Range('E38').Value = '  -2.56%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D39') '5.971'
$ws.Range('E39').Value = '  +1.92%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D40') '0.01586'
$ws.Range('E40').Value = '  -1.41%  '
$ws.Range('D41').Value = '1.057.31'
$ws.Range('E41').Value = '  +0.63%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range('D42') '104.82'
$ws.Range('E42').Value = '  +0.34%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws.Range('D43') '1.005'
$ws.Range('E43').Value = '  -0.40%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range('D44') '0.8413'
$ws.Range('E44').Value = '  -2.85%  '
$ws.Range('D45').Value = '1.788.96'
$ws.Range('E45').Value = '  -0.79%  '
Set-TextValue $ws.Range('D46') '57.42'
$ws.Range('E46').Value = '  -0.52%  '
Set-TextValue $ws.Range('D47') '1.007'
$ws.Range('E47').Value = '  -0.69%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₈104'
$ws.Range('E48').Value = '  -2.21%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D49') '0.05318'
$ws.Range('E49').Value = '  +2.95%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range('D50') '0.4340'
$ws.Range('E50').Value = '  -0.94%  '
Set-TextValue $ws.Range('D51') '7.932'
$ws.Range('E51').Value = '  -1.19%  '
